$d = $word.ActiveDocument

# "The structure were solved by" -> "The structure was solved by"
$d.Content.Find.Execute("structure were solved", $true, $false, $false, $false, $false,
                         $true, 1, $false, "structure was solved", 2)
